# Insert a new data row at row 734 (this pushes the existing rows
# 734..791 down to 735..792, extending the sheet's used range from
# A1:T791 to A1:T792), then populate the newly inserted row with the
# new "Frutilla" price-report record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(734).Insert()

$ws.Range("A734").Value = 5
$ws.Range("B734").Value = "Macroferia Regional de Talca"
$ws.Range("C734").Value = "Maule"
$ws.Range("D734").Value = 45013
$ws.Range("E734").Value = 7
$ws.Range("F734").Value = "Fruta"
$ws.Range("G734").Value = 100101
$ws.Range("H734").Value = "Berries"
$ws.Range("I734").Value = 100112025
$ws.Range("J734").Value = "Frutilla"
$ws.Range("K734").Value = "Sin especificar"
$ws.Range("L734").Value = "Primera"
$ws.Range("M734").Value = 150
$ws.Range("N734").Value = 7000
$ws.Range("O734").Value = 7000
$ws.Range("P734").Value = 7000
$ws.Range("Q734").Value = "$/caja 7 kilos"
$ws.Range("R734").Value = "Región del Maule"
$ws.Range("S734").Value = 1000
$ws.Range("T734").Value = 7
